# Add a new experiment row (row 34) describing EXP32: regularized TPR run
# with regularization weight = 0.00001, following the same pattern as the
# preceding rows (25-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting from the row above (row 33) so the new row keeps
# the same fills/borders/wrap settings without touching unrelated columns. ---
$ws.Range("A33:F33").Copy()
$ws.Range("A34:F34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K33:L33").Copy()
$ws.Range("K34:L34").PasteSpecial(-4122)   # xlPasteFormats

# --- Row height matches the extra-long description wrapping to 180pt. ---
$ws.Rows.Item(34).RowHeight = 180

# --- Column A: rich-text experiment description, mixing bold call-outs
# with normal text (mirrors the formatting used in rows 28-33). ---
$descA = "Just TPR no LSTM in `nphrase embedding layer `nbatchsize = 40. With visualizations. With regularization. Regularization weights=0.00001 [Resuming from latest successful commit, running from QA_TPR_for_Run_TPRregularizationFinal]. "
$cellA = $ws.Range("A34")
$cellA.Value = $descA

$boldRuns = @("With visualizations", "With regularization", "Regularization weights=0.00001")
$searchFrom = 0
$prevEnd = 0
foreach ($run in $boldRuns) {
    $start = $descA.IndexOf($run, $searchFrom)
    $len = $run.Length

    # normal run before this bold phrase
    if ($start -gt $prevEnd) {
        $cellA.Characters($prevEnd + 1, $start - $prevEnd).Font.Bold = $false
    }
    # the bold phrase itself
    $cellA.Characters($start + 1, $len).Font.Bold = $true

    $prevEnd = $start + $len
    $searchFrom = $prevEnd
}
# trailing normal run after the last bold phrase
if ($prevEnd -lt $descA.Length) {
    $cellA.Characters($prevEnd + 1, $descA.Length - $prevEnd).Font.Bold = $false
}

# --- Column D: logfile name ---
$ws.Range("D34").Value = "EXP32.txt"

# --- Column C: machine / GPU ---
$ws.Range("C34").Value = "DLDGX / 7"

# --- Column B: command used to launch the run ---
$ws.Range("B34").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --justTPR True --TPRregularizer1 True --TPRvis True --cF 0.00001 --cR 0.00001 --batch_size 40 --run_id 29 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP32.txt"

# --- Column E: run_id ---
$ws.Range("E34").Value = 29

# --- Column F: pane number in tmux ---
$ws.Range("F34").Value = 0

# Keep the selection on A33 (matches the pre-existing selection state).
$ws.Range("A33").Select() | Out-Null
